# Applies the "Add files via upload" edit: extends the existing practice
# workbook (Sheet1) with a couple of extra percentage-arithmetic examples
# and a "division/percentage" column block for the student pass/fail
# table, then fills in two previously-empty sheets (Sheet2: nested-IF
# grading + 5-subject percentage/division table; Sheet3: salary HRA/DA
# slab calculation) that were only stubs before.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1 -------------------------------------------------------------
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

# Small "100 * 10%" worked-example block added next to the existing
# sum/div/mul table (rows 3-9).
$ws1.Range("D3").Value = "100 KA 10 %"
$ws1.Range("E3").Value = "100*0.10"
$ws1.Range("F3").Formula = "=100*0.1"

$ws1.Range("F4").Formula = "=154333*0.1"

$ws1.Range("E6").Formula = "=100*10/100"

$ws1.Range("E8").Value = "100*10/100"
$ws1.Range("F8").Value = "100/1 * 10/100"

$ws1.Range("F9").Value = "100/1 * .10"

# Extend the pass/fail student table (rows 41-45 existed already) with
# a "total" and "percentage" column (G/H), then add two more students
# (rows 46-47) and a couple of trailing one-cell notes (rows 48-49).
$ws1.Range("G41").Formula = "=SUM(B41:D41)"
$ws1.Range("H41").Formula = "=G41*100/300"

$ws1.Range("G42:G45").Formula = "=SUM(B42:D42)"
$ws1.Range("H42:H45").Formula = "=G42*100/300"

$ws1.Range("A46").Value = "f"
$ws1.Range("B46").Value = 56
$ws1.Range("C46").Value = 34
$ws1.Range("D46").Value = 22
$ws1.Range("E46").Formula = "=COUNTIF(B46:D46,""<30"")"
$ws1.Range("F46").Formula = "=IF(E46=0,""pass"",""Faill"")"
$ws1.Range("G46").Formula = "=SUM(B46:D46)"
$ws1.Range("H46").Formula = "=G46*100/300"

$ws1.Range("A47").Value = "g"
$ws1.Range("B47").Value = 11
$ws1.Range("C47").Value = 30
$ws1.Range("D47").Value = 45
$ws1.Range("E47").Formula = "=COUNTIF(B47:D47,""<30"")"
$ws1.Range("F47").Formula = "=IF(E47=0,""pass"",""Faill"")"
$ws1.Range("G47").Formula = "=SUM(B47:D47)"
$ws1.Range("H47").Formula = "=G47*100/300"

$ws1.Range("A48").Value = "h"
$ws1.Range("A49").Value = "i"

$ws1.Range("F9").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet2 ---------------------------------------------------------------
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Range("A1").Value = "if else /nested if "
$ws2.Range("G1").Value = 0
$ws2.Range("H1").Value = "pass"
$ws2.Range("G2").Value = 1
$ws2.Range("H2").Value = "suplmentry "

$ws2.Range("A3").Value = "Student Name"
$ws2.Range("B3").Value = "hindi"
$ws2.Range("C3").Value = "math"
$ws2.Range("D3").Value = "eng"
$ws2.Range("G3").Value = "1> "
$ws2.Range("H3").Value = "fail"

$ws2.Range("A4").Value = "a"
$ws2.Range("B4").Value = 23
$ws2.Range("C4").Value = 34
$ws2.Range("D4").Value = 54
$ws2.Range("E4").Formula = "=COUNTIF(B4:D4,""<30"")"
$ws2.Range("G4").Formula = "=IF(E4=0,""pass"",IF(E4=1,""suplmentry"",""fail""))"

$ws2.Range("A5").Value = "b"
$ws2.Range("B5").Value = 34
$ws2.Range("C5").Value = 54
$ws2.Range("D5").Value = 45
$ws2.Range("E5:E10").Formula = "=COUNTIF(B5:D5,""<30"")"
$ws2.Range("G5:G10").Formula = "=IF(E5=0,""pass"",IF(E5=1,""suplmentry"",""fail""))"

$ws2.Range("A6").Value = "c"
$ws2.Range("B6").Value = 45
$ws2.Range("C6").Value = 32
$ws2.Range("D6").Value = 65

$ws2.Range("A7").Value = "d"
$ws2.Range("B7").Value = 26
$ws2.Range("C7").Value = 23
$ws2.Range("D7").Value = 76

$ws2.Range("A8").Value = "e"
$ws2.Range("B8").Value = 54
$ws2.Range("C8").Value = 56
$ws2.Range("D8").Value = 56

$ws2.Range("A9").Value = "f"
$ws2.Range("B9").Value = 56
$ws2.Range("C9").Value = 34
$ws2.Range("D9").Value = 22

$ws2.Range("A10").Value = "g"
$ws2.Range("B10").Value = 11
$ws2.Range("C10").Value = 20
$ws2.Range("D10").Value = 45

$ws2.Range("A11").Value = "h"

# Second block: 5-subject percentage/division table.
$ws2.Range("A16").Value = "Name"
$ws2.Range("B16").Value = "Hindi"
$ws2.Range("C16").Value = "Maths"
$ws2.Range("D16").Value = "English"
$ws2.Range("E16").Value = "Arts"
$ws2.Range("F16").Value = "History"
$ws2.Range("G16").Value = "Percentage "
$ws2.Range("H16").Value = "Division"

$ws2.Range("A17").Value = "A"
$ws2.Range("B17").Value = 76
$ws2.Range("C17").Value = 50
$ws2.Range("D17").Value = 60
$ws2.Range("E17").Value = 60
$ws2.Range("F17").Value = 60
$ws2.Range("G17").Formula = "=SUM(B17:F17)*100/500"
$ws2.Range("H17").Value = "If(=60 ""First Division)"" (If 50-59=""Second Division)"

$ws2.Range("A18").Value = "B"
$ws2.Range("B18").Value = 30
$ws2.Range("C18").Value = 30
$ws2.Range("D18").Value = 40
$ws2.Range("E18").Value = 40
$ws2.Range("F18").Value = 70
$ws2.Range("G18:G21").Formula = "=SUM(B18:F18)*100/500"

$ws2.Range("A19").Value = "C"
$ws2.Range("B19").Value = 20
$ws2.Range("C19").Value = 60
$ws2.Range("D19").Value = 50
$ws2.Range("E19").Value = 40
$ws2.Range("F19").Value = 30

$ws2.Range("A20").Value = "D"
$ws2.Range("B20").Value = 60
$ws2.Range("C20").Value = 70
$ws2.Range("D20").Value = 80
$ws2.Range("E20").Value = 50
$ws2.Range("F20").Value = 80

$ws2.Range("A21").Value = "E"
$ws2.Range("B21").Value = 30
$ws2.Range("C21").Value = 30
$ws2.Range("D21").Value = 40
$ws2.Range("E21").Value = 50
$ws2.Range("F21").Value = 50

$ws2.Range("E13").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet3 ---------------------------------------------------------------
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")

$ws3.Range("A1").Value = "Name "
$ws3.Range("B1").Value = "Salary"
$ws3.Range("C1").Value = "HRA"
$ws3.Range("D1").Value = "DA"
$ws3.Range("E1").NumberFormat = "0%"

$ws3.Range("A2").Value = "A"
$ws3.Range("B2").Value = 6000
$ws3.Range("C2").Formula = "=IF(B2>5000,IF(B2<10000,B2*0.1,B2*0.15),B2*0.15)"
$ws3.Range("C2").NumberFormat = "0"
$ws3.Range("E2").NumberFormat = "0%"

$ws3.Range("A3").Value = "B"
$ws3.Range("B3").Value = 7000
$ws3.Range("C3:C6").Formula = "=IF(B3>5000,IF(B3<10000,B3*0.1,B3*0.15),B3*0.15)"
$ws3.Range("C3:C6").NumberFormat = "0"
$ws3.Range("E3").NumberFormat = "0%"

$ws3.Range("A4").Value = "C"
$ws3.Range("B4").Value = 14000
$ws3.Range("E4").NumberFormat = "0%"

$ws3.Range("A5").Value = "D"
$ws3.Range("B5").Value = 5000
$ws3.Range("E5").NumberFormat = "0%"

$ws3.Range("A6").Value = "E"
$ws3.Range("B6").Value = 12000

$ws3.PageSetup.Orientation = 1

$ws3.Range("C3").Select() | Out-Null

$ws1.Activate() | Out-Null
